$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, reusing the same formatting as the
# other header cells (e.g. G1: bold, bordered, centered header style).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the data values for the new Save column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
